$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; insert it as the new first
# data row (row 58), pushing the existing rows 58:70 down to 59:71.
$ws.Rows("58").Insert()

$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 45015
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112042
$ws.Range("G58").Value = "Locoto"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 90
$ws.Range("K58").Value = 4400
$ws.Range("L58").Value = 4400
$ws.Range("M58").Value = 4400
$ws.Range("N58").Value = '$/kilo'
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 4400
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = "Hortaliza"
